# "icons on each button vers2"
# Add a new worksheet named "icons" after the last existing sheet
# ("comments"), make it the active sheet/tab, and populate it with the
# "hdd.png" icon filename used for the OS-Drive / Data-Drive buttons
# (mirrors the layout already used on the "comments" sheet).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet so it lands at
# the end of the tab strip (sheetId 5 / 5th tab).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$icons = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$icons.Name = "icons"

# Make it the active sheet -> bumps workbook activeTab to 4 and moves
# tabSelected onto this sheet.
$icons.Activate()

# Icon filename shown for both drive buttons.
$icons.Range("A2").Value = "hdd.png"
$icons.Range("B2").Value = "hdd.png"

# Apply the sheet's normal/default cell style across the small
# "button" block (same staircase footprint as authored originally).
$icons.Range("A2:D3").Style = "Normal"
$icons.Range("B4").Style = "Normal"
$icons.Range("D4").Style = "Normal"
$icons.Range("B5").Style = "Normal"
$icons.Range("D5").Style = "Normal"
$icons.Range("D6").Style = "Normal"

# Leave the cursor where the sheet was saved with it selected.
$icons.Range("D16").Select() | Out-Null
